$wb = $excel.ActiveWorkbook

# Add the new sheet after the last existing sheet ("companion") so it lands
# at the end of the tab order, then rename it to "companion2H".
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "companion2H"

# Header row (matches the other sheets: Date / Chapter).
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Chapter"

# Data rows for the companion2H group (dates as Excel serials, matching the
# other sheets, 2023-06-07 / 2023-06-13 / 2023-06-07).
$ws.Range("A2").Value = 45084
$ws.Range("B2").Value = "Chapter 17"

$ws.Range("A3").Value = 45090
$ws.Range("B3").Value = "Chapter 18"

$ws.Range("A4").Value = 45084
$ws.Range("B4").Value = "Chapter 19"

# Match the date formatting used on the sibling sheets.
$ws.Range("A2:A4").NumberFormat = "yyyy-mm-dd"
